# modul receiving update 10/24/25
# Updates part numbers in the Marcone cart-vs-inventory table and appends
# the newly-received line items as additional rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Update existing part-number cells (column 1) in place ---------------
$t.Cell(2, 1).Range.Text = "5304524714"
$t.Cell(3, 1).Range.Text = "5304532207"
$t.Cell(4, 1).Range.Text = "154773201"
$t.Cell(5, 1).Range.Text = "DC97-16350U"
$t.Cell(6, 1).Range.Text = "6600JB3001C"
$t.Cell(7, 1).Range.Text = "240599803"

# --- Append new rows for the additional received items -------------------
# Each entry: PartNumber, QtyInCart, InStock?, QtyInStock, Location, Name
$newRows = @(
    @("240579820",   "1", "No", "0", "-", "-"),
    @("5303918344",  "1", "No", "0", "-", "-"),
    @("WR55X11070",  "1", "No", "0", "-", "-"),
    @("W10807577",   "1", "No", "0", "-", "-"),
    @("W10807577EXCR","1","No", "0", "-", "-"),
    @("W10859573",   "1", "No", "0", "-", "-"),
    @("242193212",   "1", "No", "0", "-", "-"),
    @("MDS65210402", "1", "No", "0", "-", "-"),
    @("242219206",   "1", "No", "0", "-", "-"),
    @("154579101",   "1", "No", "0", "-", "-"),
    @("154756401",   "2", "No", "0", "-", "-"),
    @("WB02X10400",  "4", "No", "0", "-", "-")
)

foreach ($rowData in $newRows) {
    $newRow = $t.Rows.Add()
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $newRow.Cells.Item($c).Range.Text = $rowData[$c - 1]
    }
}

Write-Host "Updated part numbers and appended" $newRows.Length "new rows."
